$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was recorded for this product/category
# combination. Insert a new row at 809 (pushing the existing rows 809-848
# down to 810-849) and populate it with the new reading.
$ws.Rows.Item(809).Insert()

$ws.Cells.Item(809, 1).Value = 8
$ws.Cells.Item(809, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(809, 3).Value = "Coquimbo"
$ws.Cells.Item(809, 4).Value = 44939
$ws.Cells.Item(809, 5).Value = 4
$ws.Cells.Item(809, 6).Value = 100112045
$ws.Cells.Item(809, 7).Value = "Zapallo"
$ws.Cells.Item(809, 8).Value = "Camote"
$ws.Cells.Item(809, 9).Value = "1a (cosecha)"
$ws.Cells.Item(809, 10).Value = 1600
$ws.Cells.Item(809, 11).Value = 800
$ws.Cells.Item(809, 12).Value = 900
$ws.Cells.Item(809, 13).Value = 850
$ws.Cells.Item(809, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(809, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(809, 16).Value = 850
$ws.Cells.Item(809, 17).Value = 1
$ws.Cells.Item(809, 18).Value = "Hortaliza"
